{"js": "// Updated units from CPUh to Core-h (and a couple of small copy fixes\n// that were bundled in the same commit).\n\nconst body = context.document.body;\n\n// 1) \"CPUh\" -> \"Core-h\" everywhere it appears (3 occurrences: the\n//    \"Total CPUh\" label, the cyan-highlighted instruction text, and the\n//    \"total CPUh\" justification bullet).\nconst cpuh = body.search(\"CPUh\", { matchCase: true, matchWholeWord: false });\ncpuh.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < cpuh.items.length; i++) {\n  cpuh.items[i].insertText(\"Core-h\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Drop the stray \"Early \" before \"Access\" and lowercase \"Access\" so\n//    the sentence reads \"...applying for Cirrus access. This should cover:\"\nconst early = body.search(\"Early \", { matchCase: true });\nearly.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < early.items.length; i++) {\n  early.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst access = body.search(\"Access. This should cover:\", { matchCase: true });\naccess.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < access.items.length; i++) {\n  access.items[i].insertText(\"access. This should cover:\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Updated units from CPUh to Core-h (and a couple of small copy fixes\n# that were bundled in the same commit).\n\n$d = $word.ActiveDocument\n\n# 1) \"CPUh\" -> \"Core-h\" everywhere it appears (3 occurrences: the\n#    \"Total CPUh\" label, the cyan-highlighted instruction text, and the\n#    \"total CPUh\" justification bullet).\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Text = \"CPUh\"\n$rng1.Find.Replacement.Text = \"Core-h\"\n$rng1.Find.Forward = $true\n$rng1.Find.Wrap = 1\n$rng1.Find.Execute($rng1.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng1.Find.Replacement.Text, 2)\n\n# 2) Drop the stray \"Early \" before \"Access\" and lowercase \"Access\" so\n#    the sentence reads \"...applying for Cirrus access. This should cover:\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Text = \"Early \"\n$rng2.Find.Replacement.Text = \"\"\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 1\n$rng2.Find.Execute($rng2.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng2.Find.Replacement.Text, 2)\n\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$rng3.Find.Text = \"Access. This should cover:\"\n$rng3.Find.Replacement.Text = \"access. This should cover:\"\n$rng3.Find.Forward = $true\n$rng3.Find.Wrap = 1\n$rng3.Find.Execute($rng3.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng3.Find.Replacement.Text, 2)\n"}
